{"js": "// Modify protocol field in school transports' templates.\n// The run that reads \".: ${\" (immediately before the \"protocol\" merge-field\n// run) needs to become \".: \u03a6.15.1/${\" so the finished sentence renders as\n// \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: \u03a6.15.1/${protocol}\".\nconst body = context.document.body;\n\nconst results = body.search(\".: ${\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the target text \".: ${\" in the document.');\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\".: \u03a6.15.1/${\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Modify protocol field in school transports' templates.\n# The run that reads \".: ${\" (immediately before the \"protocol\" merge-field\n# run) needs to become \".: \u03a6.15.1/${\" so the finished sentence renders as\n# \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: \u03a6.15.1/${protocol}\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '.: ${'\n$find.Replacement.Text = '.: \u03a6.15.1/${'\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
